{"js": "// Locate the \"Define the default value of expense_id of ExpenseDetail\n// model\" bullet and replace it with a run-split version (spell-check\n// proofErr markers around expense_id / ExpenseDetail) followed by three\n// new list bullets, the middle one being the 'tag' field removal bullet\n// from the commit message. Search by text so this is resilient to the\n// paragraph not being last; fall back to the last paragraph (its usual\n// position) if the text isn't found verbatim.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Define the default value of expense_id of ExpenseDetail model\";\nlet anchorParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf(anchorText) === 0) {\n    anchorParagraph = p;\n  }\n}\nif (!anchorParagraph) {\n  anchorParagraph = paragraphs.items[paragraphs.items.length - 1];\n}\n\nconst targetRange = anchorParagraph.getRange(\"Whole\");\n\nconst ooxml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t xml:space=\"preserve\">Define the default value of </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>expense_id</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> of </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>ExpenseDetail</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> model</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t xml:space=\"preserve\">Remove all id fields from </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>attr_accessible</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t>Remove tag field from companies table</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t xml:space=\"preserve\">Remove </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>CompanyId</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> from Expenses table</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Define the default value of expense_id of ExpenseDetail\n# model\" bullet - the new bullets get appended right after it. Search by\n# its text so the script is resilient to the paragraph not being last;\n# fall back to the last paragraph of the document (its usual position).\n$anchorPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Define the default value of expense_id of ExpenseDetail model*\") {\n        $anchorPara = $p\n    }\n}\nif ($anchorPara -eq $null) {\n    $anchorPara = $d.Paragraphs.Last\n}\n$target = $d.Range($anchorPara.Range.Start, $anchorPara.Range.End)\n\n# Replace that paragraph's content with a version whose run is split\n# around the spell-check-flagged tokens (expense_id / ExpenseDetail),\n# and append three new list bullets after it - the middle one being the\n# 'tag' field removal bullet called out in the commit message.\n$xml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t xml:space=\"preserve\">Define the default value of </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>expense_id</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> of </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>ExpenseDetail</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> model</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t xml:space=\"preserve\">Remove all id fields from </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>attr_accessible</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t>Remove tag field from companies table</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r><w:t xml:space=\"preserve\">Remove </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>CompanyId</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> from Expenses table</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$target.InsertXML($xml)\n"}
